# Applies the "Add files via upload" edit to Sheet1 of the workbook:
#  - Fill in the running S.No. (column A) for rows 348..400 (values 347..399)
#  - Fix a typo: D351 should reference "haseen" instead of "hasen"
#    (the now-unused "hasen" shared string is pruned automatically on save)
#  - Update the sheet view's scroll position / active selection to D352

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A numbering continues for rows 348 through 400 (value = row - 1).
for ($r = 348; $r -le 400; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Correct the name in D351 from "hasen" to "haseen".
$ws.Range("D351").Value = "haseen"

# Move the view: scroll so row 339 is at the top and select D352.
$ws.Range("D352").Select()
$excel.ActiveWindow.ScrollRow = 339
$excel.ActiveWindow.ScrollColumn = 1
